$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.5019577983990334
$ws.Range("C2").Value = 0.5019577983990334
$ws.Range("D2").Value = 0.3105127224087574
$ws.Range("E2").Value = 0.5572366843709748
$ws.Range("F2").Value = 0.251107588723997
$ws.Range("G2").Value = 14

$ws.Range("B3").Value = 0.3307979378781269
$ws.Range("C3").Value = 0.3343319746521877
$ws.Range("D3").Value = 0.1645692088094062
$ws.Range("E3").Value = 0.4056713063668741
$ws.Range("F3").Value = 0.244411730890862
$ws.Range("G3").Value = 13

$ws.Range("B4").Value = 0.2377396291553046
$ws.Range("C4").Value = 0.2452017974494533
$ws.Range("D4").Value = 0.09171268429751143
$ws.Range("E4").Value = 0.3028410214906683
$ws.Range("F4").Value = 0.1959384495933066
$ws.Range("G4").Value = 12

$ws.Range("B5").Value = 0.373625158811909
$ws.Range("C5").Value = 0.373625158811909
$ws.Range("D5").Value = 0.1693451806315346
$ws.Range("E5").Value = 0.4115157112815191
$ws.Range("F5").Value = 0.1808987657993869
$ws.Range("G5").Value = 11

$ws.Range("B6").Value = 0.3392289017490621
$ws.Range("C6").Value = 0.3392289017490621
$ws.Range("D6").Value = 0.1435050013861421
$ws.Range("E6").Value = 0.3788205398155466
$ws.Range("F6").Value = 0.1777287371382058
$ws.Range("G6").Value = 10

$ws.Range("B7").Value = 0.323658031409521
$ws.Range("C7").Value = 0.323658031409521
$ws.Range("D7").Value = 0.1411500930911699
$ws.Range("E7").Value = 0.3756994717738766
$ws.Range("F7").Value = 0.202348754060147
$ws.Range("G7").Value = 9

$ws.Range("B8").Value = 0.3534922319734252
$ws.Range("C8").Value = 0.3534922319734252
$ws.Range("D8").Value = 0.1628866010887078
$ws.Range("E8").Value = 0.4035921221836568
$ws.Range("F8").Value = 0.2082027666157719
$ws.Range("G8").Value = 8

$ws.Range("B9").Value = 0.3545941574938994
$ws.Range("C9").Value = 0.3545941574938994
$ws.Range("D9").Value = 0.1673332494109491
$ws.Range("E9").Value = 0.4090638695985617
$ws.Range("F9").Value = 0.2202928468255327
$ws.Range("G9").Value = 7

$ws.Range("B10").Value = 0.3114747142229594
$ws.Range("C10").Value = 0.3114747142229594
$ws.Range("D10").Value = 0.1227287858089793
$ws.Range("E10").Value = 0.3503266844089661
$ws.Range("F10").Value = 0.1756551902177846
$ws.Range("G10").Value = 6

$ws.Range("B11").Value = 0.3600284238150658
$ws.Range("C11").Value = 0.3600284238150658
$ws.Range("D11").Value = 0.1554279835465356
$ws.Range("E11").Value = 0.3942435586620733
$ws.Range("F11").Value = 0.1796090114379531
$ws.Range("G11").Value = 5
